$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").ClearContents()
$ws.Range("N26").Value = 0
$ws.Range("H116").Value = 3581958.2
$ws.Range("I116").Value = 14288354
$ws.Range("J116").Value = 13159.934
$ws.Range("K116").Value = 14288354
$ws.Range("L116").Value = 13159.934
$ws.Range("M116").Value = -14284912
$ws.Range("N116").Value = -20043.934
$ws.Range("H121").Value = 1026.25
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1026.25
$ws.Range("K121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("M121").Value = 3078.75
$ws.Range("N121").Value = -6572.75
$ws.Range("H137").Value = 1051.6111
$ws.Range("I137").Value = 861.2308
$ws.Range("J137").Value = 1546.6
$ws.Range("K137").Value = 2583.6924
$ws.Range("L137").Value = 4639.799999999999
$ws.Range("M137").Value = -33.69239999999991
$ws.Range("N137").Value = -9739.799999999999
$ws.Range("H138").Value = 3838.1
$ws.Range("I138").Value = 869.119
$ws.Range("J138").Value = 7119.6055
$ws.Range("K138").Value = 2607.357
$ws.Range("L138").Value = 21358.8165
$ws.Range("M138").Value = 2532.643
$ws.Range("N138").Value = -31638.8165

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 6000
$ws.Range("J11").Value = 6000
$ws.Range("L11").Value = 6000
$ws.Range("N11").Value = -6288
$ws.Range("H32").Value = 8111.488
$ws.Range("I32").Value = 4056.2058
$ws.Range("J32").Value = 27808.572
$ws.Range("K32").Value = 4056.2058
$ws.Range("L32").Value = 27808.572
$ws.Range("M32").Value = -3769.2058
$ws.Range("N32").Value = -28382.572
$ws.Range("H61").Value = 1414.3158
$ws.Range("I61").Value = 1177.6
$ws.Range("J61").Value = 2302
$ws.Range("K61").Value = 1177.6
$ws.Range("L61").Value = 2302
$ws.Range("M61").Value = -965.5999999999999
$ws.Range("N61").Value = -2726
$ws.Range("H74").Value = 947.9666999999999
$ws.Range("I74").Value = 872.4583
$ws.Range("J74").Value = 1250
$ws.Range("K74").Value = 872.4583
$ws.Range("L74").Value = 1250
$ws.Range("M74").Value = 1.541699999999992
$ws.Range("N74").Value = -2998
$ws.Range("H76").Value = 36144
$ws.Range("J76").Value = 36144
$ws.Range("L76").Value = 36144
$ws.Range("N76").Value = -36820
$ws.Range("H77").Value = 947.9666999999999
$ws.Range("I77").Value = 872.4583
$ws.Range("J77").Value = 1250
$ws.Range("K77").Value = 4362.2915
$ws.Range("L77").Value = 6250
$ws.Range("M77").Value = 5.708499999999731
$ws.Range("N77").Value = -14986
$ws.Range("H79").Value = 36144
$ws.Range("J79").Value = 36144
$ws.Range("L79").Value = 36144
$ws.Range("N79").Value = -38484
$ws.Range("H101").Value = 37290.727
$ws.Range("J101").Value = 37290.727
$ws.Range("L101").Value = 37290.727
$ws.Range("N101").Value = -43780.727
$ws.Range("H132").Value = 2170.1936
$ws.Range("I132").Value = 1608.8889
$ws.Range("J132").Value = 2947.3845
$ws.Range("K132").Value = 4826.6667
$ws.Range("L132").Value = 8842.1535
$ws.Range("M132").Value = -2296.6667
$ws.Range("N132").Value = -13902.1535
$ws.Range("H136").Value = 1414.3158
$ws.Range("I136").Value = 1177.6
$ws.Range("J136").Value = 2302
$ws.Range("K136").Value = 3532.8
$ws.Range("L136").Value = 6906
$ws.Range("M136").Value = -982.7999999999997
$ws.Range("N136").Value = -12006

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1411.1277
$ws.Range("I134").Value = 953.71875
$ws.Range("J134").Value = 2386.9333
$ws.Range("K134").Value = 2861.15625
$ws.Range("L134").Value = 7160.7999
$ws.Range("M134").Value = -326.15625
$ws.Range("N134").Value = -12230.7999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").ClearContents()
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = 0
$ws.Range("H58").Value = 1087.5
$ws.Range("I58").Value = 818.2941
$ws.Range("J58").Value = 2002.8
$ws.Range("K58").Value = 818.2941
$ws.Range("L58").Value = 2002.8
$ws.Range("M58").Value = -615.2941
$ws.Range("N58").Value = -2408.8
$ws.Range("H134").Value = 1266.975
$ws.Range("I134").Value = 1101.9395
$ws.Range("J134").Value = 2045
$ws.Range("K134").Value = 3305.8185
$ws.Range("L134").Value = 6135
$ws.Range("M134").Value = -770.8184999999999
$ws.Range("N134").Value = -11205
$ws.Range("H136").Value = 1087.5
$ws.Range("I136").Value = 818.2941
$ws.Range("J136").Value = 2002.8
$ws.Range("K136").Value = 2454.8823
$ws.Range("L136").Value = 6008.4
$ws.Range("M136").Value = 95.11770000000024
$ws.Range("N136").Value = -11108.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 984.80646
$ws.Range("I5").Value = 335.31033
$ws.Range("J5").Value = 10402.5
$ws.Range("K5").Value = 1005.93099
$ws.Range("L5").Value = 31207.5
$ws.Range("M5").Value = -893.9309900000001
$ws.Range("N5").Value = -31431.5
$ws.Range("H33").Value = 94.5
$ws.Range("I33").Value = 89
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 534
$ws.Range("L33").Value = 600
$ws.Range("M33").Value = -251
$ws.Range("N33").Value = -1166
$ws.Range("H64").Value = 3500
$ws.Range("I64").Value = 1000
$ws.Range("J64").Value = 4750
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 14250
$ws.Range("M64").Value = -2730
$ws.Range("N64").Value = -14790
$ws.Range("H67").Value = 3500
$ws.Range("I67").Value = 1000
$ws.Range("J67").Value = 4750
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 14250
$ws.Range("M67").Value = -2064
$ws.Range("N67").Value = -16122
$ws.Range("H113").Value = 2551394.5
$ws.Range("I113").Value = 363.85367
$ws.Range("J113").Value = 9524212
$ws.Range("K113").Value = 1091.56101
$ws.Range("L113").Value = 28572636
$ws.Range("M113").Value = 1078.43899
$ws.Range("N113").Value = -28576976
$ws.Range("H131").Value = 903.17145
$ws.Range("J131").Value = 996.5714
$ws.Range("L131").Value = 2989.7142
$ws.Range("N131").Value = -13069.7142
$ws.Range("H135").Value = 984.80646
$ws.Range("I135").Value = 335.31033
$ws.Range("J135").Value = 10402.5
$ws.Range("K135").Value = 3017.79297
$ws.Range("L135").Value = 93622.5
$ws.Range("M135").Value = -482.7929700000004
$ws.Range("N135").Value = -98692.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 4426.75
$ws.Range("I3").Value = 3203
$ws.Range("J3").Value = 4834.6665
$ws.Range("K3").Value = 3203
$ws.Range("L3").Value = 4834.6665
$ws.Range("M3").Value = -3087
$ws.Range("N3").Value = -5066.6665
$ws.Range("H132").Value = 6221.1665
$ws.Range("I132").Value = 7730.375
$ws.Range("J132").Value = 3202.75
$ws.Range("K132").Value = 23191.125
$ws.Range("L132").Value = 9608.25
$ws.Range("M132").Value = -20661.125
$ws.Range("N132").Value = -14668.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").ClearContents()
$ws.Range("N87").Value = 0
$ws.Range("H88").Value = 10000
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 10000
$ws.Range("L88").ClearContents()
$ws.Range("M88").Value = -9572
$ws.Range("N88").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").ClearContents()
$ws.Range("N90").Value = 0
$ws.Range("H91").Value = 10000
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 10000
$ws.Range("L91").ClearContents()
$ws.Range("M91").Value = -8518
$ws.Range("N91").Value = 0
$ws.Range("H104").Value = 29407.7
$ws.Range("J104").Value = 29407.7
$ws.Range("L104").Value = 29407.7
$ws.Range("N104").Value = -36395.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 22843.625
$ws.Range("J63").Value = 25178.428
$ws.Range("L63").Value = 25178.428
$ws.Range("N63").Value = -26426.428
$ws.Range("H66").Value = 22843.625
$ws.Range("J66").Value = 25178.428
$ws.Range("L66").Value = 75535.284
$ws.Range("N66").Value = -81775.284
$ws.Range("H95").Value = 38133.332
$ws.Range("J95").Value = 37560
$ws.Range("L95").Value = 37560
$ws.Range("N95").Value = -43052
$ws.Range("H103").Value = 25390.4
$ws.Range("J103").Value = 25390.4
$ws.Range("L103").Value = 25390.4
$ws.Range("N103").Value = -27734.4
$ws.Range("H105").Value = 33514.75
$ws.Range("J105").Value = 33514.75
$ws.Range("L105").Value = 33514.75
$ws.Range("N105").Value = -40502.75
$ws.Range("H132").Value = 1079.8158
$ws.Range("I132").Value = 750.7308
$ws.Range("J132").Value = 1792.8334
$ws.Range("K132").Value = 2252.1924
$ws.Range("L132").Value = 5378.5002
$ws.Range("M132").Value = 277.8076000000001
$ws.Range("N132").Value = -10438.5002
